$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("L9").Value2 = -0.3802
$ws.Range("M9").Value2 = -0.1156
$ws.Range("N9").Value2 = -0.0439
$ws.Range("O9").Value2 = -0.2949
$ws.Range("P9").Value2 = -0.4257
$ws.Range("Q9").Value2 = -0.0613
$ws.Range("R9").Value2 = -0.0515
$ws.Range("S9").Value2 = -0.3876
$ws.Range("T9").Value2 = -32.9216

# Row 13
$ws.Range("L13").Value2 = 0.0179
$ws.Range("M13").Value2 = 0.0029
$ws.Range("N13").Value2 = -0.0022
$ws.Range("P13").Value2 = -0.0113
$ws.Range("Q13").Value2 = 0.007
$ws.Range("R13").Value2 = -0.0111
$ws.Range("S13").Value2 = -0.0244
$ws.Range("T13").Value2 = -5.4036

# Row 15
$ws.Range("L15").Value2 = 0.3134
$ws.Range("M15").Value2 = 0.2152
$ws.Range("N15").Value2 = 0.2641
$ws.Range("O15").Value2 = 0.5368
$ws.Range("P15").Value2 = 0.5283
$ws.Range("Q15").Value2 = 0.5569
$ws.Range("R15").Value2 = 0.5092
$ws.Range("S15").Value2 = 0.2708
$ws.Range("T15").Value2 = 7.6505

# Row 23
$ws.Range("L23").Value2 = 0.225
$ws.Range("M23").Value2 = 0.2688
$ws.Range("N23").Value2 = -0.0137
$ws.Range("O23").Value2 = -0.1003
$ws.Range("P23").Value2 = -0.0852
$ws.Range("Q23").Value2 = -0.0702
$ws.Range("R23").Value2 = -0.0351
$ws.Range("S23").Value2 = -0.0729
$ws.Range("T23").Value2 = -6.6924

# Row 31
$ws.Range("L31").Value2 = -1.2456
$ws.Range("M31").Value2 = -0.1857
$ws.Range("N31").Value2 = -0.1632
$ws.Range("O31").Value2 = -0.3603
$ws.Range("P31").Value2 = 0.0317
$ws.Range("Q31").Value2 = 0.0696
$ws.Range("R31").Value2 = 0.1879
$ws.Range("S31").Value2 = -0.0504
$ws.Range("T31").Value2 = -72.2455

# Row 35
$ws.Range("L35").Value2 = -0.1144
$ws.Range("M35").Value2 = -0.1186
$ws.Range("N35").Value2 = -0.0684
$ws.Range("O35").Value2 = -0.2598
$ws.Range("P35").Value2 = -0.3725
$ws.Range("Q35").Value2 = -0.0132
$ws.Range("R35").Value2 = 0.0133
$ws.Range("S35").Value2 = -0.3159
$ws.Range("T35").Value2 = -6.8434

# Row 43
$ws.Range("L43").Value2 = 0.0681
$ws.Range("M43").Value2 = 0.123
$ws.Range("N43").Value2 = -0.0103
$ws.Range("O43").Value2 = 0.2002
$ws.Range("P43").Value2 = 0.2917
$ws.Range("Q43").Value2 = -0.0699
$ws.Range("R43").Value2 = -0.0998
$ws.Range("S43").Value2 = 0.2271
$ws.Range("T43").Value2 = -37.433

# Row 47
$ws.Range("L47").Value2 = 0.0149
$ws.Range("M47").Value2 = 0.1232
$ws.Range("N47").Value2 = 0.0766
$ws.Range("P47").Value2 = -0.0182
$ws.Range("Q47").Value2 = -0.0672
$ws.Range("R47").Value2 = -0.0456
$ws.Range("S47").Value2 = -0.036
$ws.Range("T47").Value2 = -0.7106

# Row 57
$ws.Range("L57").Value2 = -0.0131
$ws.Range("M57").Value2 = 0.006
$ws.Range("N57").Value2 = 0.0207
$ws.Range("O57").Value2 = 0.0219
$ws.Range("P57").Value2 = 0.0242
$ws.Range("Q57").Value2 = 0.0184
$ws.Range("R57").Value2 = 0.0086
$ws.Range("S57").Value2 = 0.0021
$ws.Range("T57").Value2 = -0.0885

# Row 65
$ws.Range("L65").Value2 = -0.33
$ws.Range("M65").Value2 = -0.0631
$ws.Range("N65").Value2 = 0.0005
$ws.Range("O65").Value2 = 0.0012
$ws.Range("P65").Value2 = 0.0022
$ws.Range("Q65").Value2 = 0.0005
$ws.Range("R65").Value2 = -0.0011
$ws.Range("S65").Value2 = 0.0029
$ws.Range("T65").Value2 = -0.0709

# Row 69
$ws.Range("L69").Value2 = 0.0163
$ws.Range("M69").Value2 = 0.0449
$ws.Range("N69").Value2 = 0.0392
$ws.Range("O69").Value2 = -0.0216
$ws.Range("P69").Value2 = 0.0127
$ws.Range("Q69").Value2 = -0.0197
$ws.Range("R69").Value2 = -0.042
$ws.Range("S69").Value2 = -0.0654
$ws.Range("T69").Value2 = -0.0003

# Row 71
$ws.Range("L71").Value2 = 0.3056
$ws.Range("M71").Value2 = 0.1527
$ws.Range("N71").Value2 = 0.139
$ws.Range("O71").Value2 = 0.1345
$ws.Range("P71").Value2 = 0.1312
$ws.Range("Q71").Value2 = 0.1275
$ws.Range("R71").Value2 = 0.127
$ws.Range("S71").Value2 = -0.0085
$ws.Range("T71").Value2 = -0.3384

# Row 79
$ws.Range("L79").Value2 = 0.0992
$ws.Range("M79").Value2 = 0.0969
$ws.Range("N79").Value2 = 0.0967
$ws.Range("O79").Value2 = -0.0046
$ws.Range("P79").Value2 = -0.0044
$ws.Range("Q79").Value2 = -0.0042
$ws.Range("R79").Value2 = -0.0042
$ws.Range("S79").Value2 = -0.0037
$ws.Range("T79").Value2 = -0.0608

# Row 87
$ws.Range("L87").Value2 = 0.3152
$ws.Range("M87").Value2 = 0.2717
$ws.Range("N87").Value2 = 0.3192
$ws.Range("O87").Value2 = 0.182
$ws.Range("P87").Value2 = 0.1856
$ws.Range("Q87").Value2 = 0.1827
$ws.Range("R87").Value2 = 0.1593
$ws.Range("S87").Value2 = -0.0568
$ws.Range("T87").Value2 = -0.1023

# Row 91
$ws.Range("L91").Value2 = -0.083
$ws.Range("M91").Value2 = 0.0003
$ws.Range("N91").Value2 = 0
$ws.Range("O91").Value2 = 0.0014
$ws.Range("P91").Value2 = 0.0023
$ws.Range("Q91").Value2 = -0.0001
$ws.Range("R91").Value2 = -0.0003
$ws.Range("S91").Value2 = 0.0019
$ws.Range("T91").Value2 = -0.1259

# Row 99
$ws.Range("L99").Value2 = 0.0834
$ws.Range("M99").Value2 = -0.0003
$ws.Range("N99").Value2 = 0.0007
$ws.Range("O99").Value2 = -0.0009
$ws.Range("P99").Value2 = -0.0017
$ws.Range("Q99").Value2 = 0.0008
$ws.Range("R99").Value2 = 0.001
$ws.Range("S99").Value2 = -0.0012
$ws.Range("T99").Value2 = 0.4949

# Row 103
$ws.Range("L103").Value2 = 0.0093
$ws.Range("M103").Value2 = -0.0267
$ws.Range("N103").Value2 = -0.0272
$ws.Range("P103").Value2 = -0.0369
$ws.Range("Q103").Value2 = -0.0032
$ws.Range("R103").Value2 = 0.0052
$ws.Range("S103").Value2 = 0.0132
$ws.Range("T103").Value2 = 0.0324

# Row 113
$ws.Range("L113").Value2 = -0.0069
$ws.Range("M113").Value2 = -0.001
$ws.Range("N113").Value2 = 0.0047
$ws.Range("O113").Value2 = 0.0096
$ws.Range("P113").Value2 = 0.0167
$ws.Range("Q113").Value2 = 0.0159
$ws.Range("R113").Value2 = 0.0102
$ws.Range("S113").Value2 = 0.0061
$ws.Range("T113").Value2 = -0.009
